$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text format before writing values so that
# numeric-looking strings (e.g. "598.89", "1.00", "0.0000116") are preserved as
# literal text instead of being auto-converted to numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '65.452.34'
$ws.Range('D3').Value = '3.562.17'
$ws.Range('E3').Value = '  +3.35%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '598.89'
$ws.Range('E5').Value = '  +0.54%  '
$ws.Range('D6').Value = '140.19'
$ws.Range('E6').Value = '  +2.75%  '
$ws.Range('D7').Value = '3.562.38'
$ws.Range('E7').Value = '  +3.40%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('E9').Value = '  +0.35%  '
$ws.Range('E10').Value = '  +3.61%  '
$ws.Range('D11').Value = '7.10'
$ws.Range('E11').Value = '  -5.95%  '
$ws.Range('E12').Value = '  +3.83%  '
$ws.Range('D13').Value = '4.167.92'
$ws.Range('E13').Value = '  +3.36%  '
$ws.Range('E14').Value = '  +4.15%  '
$ws.Range('D15').Value = '3.569.29'
$ws.Range('E15').Value = '  +2.84%  '
$ws.Range('E16').Value = '  +2.36%  '
$ws.Range('E17').Value = '  +1.53%  '
$ws.Range('D18').Value = '65.394.80'
$ws.Range('E18').Value = '  -0.26%  '
$ws.Range('E19').Value = '  +3.74%  '
$ws.Range('E20').Value = '  +1.60%  '
$ws.Range('E21').Value = '  +4.02%  '
$ws.Range('D22').Value = '397.31'
$ws.Range('E22').Value = '  +0.60%  '
$ws.Range('E23').Value = '  +4.53%  '
$ws.Range('D24').Value = '3.706.52'
$ws.Range('E24').Value = '  +3.15%  '
$ws.Range('D25').Value = '74.66'
$ws.Range('E25').Value = '  +1.70%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').Value = '0.0000116'
$ws.Range('E27').Value = '  +10.53%  '
$ws.Range('D28').Value = '7.87'
$ws.Range('E28').Value = '  +7.90%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('E30').Value = '  +0.71%  '
$ws.Range('D32').Value = '3.579.80'
$ws.Range('E32').Value = '  +3.65%  '
$ws.Range('D33').Value = '23.96'
$ws.Range('E33').Value = '  +4.79%  '
$ws.Range('E35').Value = '  +1.11%  '
$ws.Range('E36').Value = '  +3.74%  '
$ws.Range('D37').Value = '7.08'
$ws.Range('E37').Value = '  +2.59%  '
$ws.Range('D38').Value = '168.92'
$ws.Range('E38').Value = '  -2.60%  '
$ws.Range('D39').Value = '1.56'
$ws.Range('E39').Value = '  +2.34%  '
$ws.Range('E40').Value = '  +4.42%  '
$ws.Range('D41').Value = '0.0803'
$ws.Range('E41').Value = '  +2.90%  '
$ws.Range('E42').Value = '  +2.13%  '
$ws.Range('D43').Value = '26.72'
$ws.Range('E43').Value = '  +16.14%  '
$ws.Range('D44').Value = '42.91'
$ws.Range('E44').Value = '  -1.48%  '
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').Value = '1.71'
$ws.Range('E46').Value = '  +4.36%  '
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').Value = '4.44'
$ws.Range('E47').Value = '  -0.19%  '
$ws.Range('D48').Value = '1.19'
$ws.Range('E48').Value = '  +7.89%  '
$ws.Range('D49').Value = '2.449.30'
$ws.Range('E49').Value = '  +10.81%  '
$ws.Range('E50').Value = '  +4.24%  '
$ws.Range('D51').Value = '2.14'
$ws.Range('E51').Value = '  +1.73%  '

# Restore the original (default) cell style now that the text values are set,
# so the cells end up with no explicit style index, matching the source format.
$dataRange.Style = "Normal"

